$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Site data: new display name (old name + " Boat Launch") and the Google Maps
# hyperlink target that used to live in column E and now lives in column F.
$sites = @(
    @{ Row = 2;  Name = "Farragut Park Boat Launch";   Url = "https://www.google.com/maps/place/Farragut+Boat+Launch/@47.9653555,-116.8259194,10z/data=!4m9!1m2!2m1!1sfarragut+park+boat+launch!3m5!1s0x53619c0ce8a3a54f:0x712d271b0454132c!8m2!3d47.9653555!4d-116.545768!15sChlmYXJyYWd1dCBwYXJrIGJvYXQgbGF1bmNoWjYKGWZhcnJ" },
    @{ Row = 3;  Name = "Pringle Park Boat Launch";    Url = "https://www.google.com/maps/place/Pringle+Park/@47.6305848,-125.5025766,5z/data=!4m9!1m2!2m1!1springle+park!3m5!1s0x0:0xa00d98ef35987353!8m2!3d48.2394254!4d-116.2936172!15sCgxwcmluZ2xlIHBhcmuSAQRwYXJr" },
    @{ Row = 4;  Name = "Trestle Creek Boat Launch";   Url = "https://www.google.com/maps/place/48%C2%B016'36.3%22N+116%C2%B020'49.5%22W/@48.2788982,-116.3555405,2006m/data=!3m1!1e3!4m5!3m4!1s0x0:0x0!8m2!3d48.276753!4d-116.347077" },
    @{ Row = 5;  Name = "Hawkins Point Boat Launch";   Url = "https://www.google.com/maps/place/Hawkins+Point/@48.2830355,-116.413729,8625m/data=!3m2!1e3!4b1!4m5!3m4!1s0x5363d8fd1c302e37:0x5d14100de1165628!8m2!3d48.2829842!4d-116.3787954" },
    @{ Row = 6;  Name = "City Beach Boat Launch";      Url = "https://www.google.com/maps/place/Sandpoint+City+Beach+Park/@48.2729824,-116.5434514,539m/data=!3m2!1e3!4b1!4m5!3m4!1s0x5363d6c49683c93d:0x59b6991355e517c7!8m2!3d48.2729789!4d-116.5412681" },
    @{ Row = 7;  Name = "Memorial Field Boat Launch";  Url = "https://www.google.com/maps/place/Lakeview+Park+Boat+Ramp/@48.2654406,-116.5592122,539m/data=!3m1!1e3!4m12!1m6!3m5!1s0x5363d6bc3607d97d:0x7de15e1319d71b43!2sWar+Memorial+Field!8m2!3d48.265787!4d-116.5597647!3m4!1s0x5363d6bdc5087dad:0x467cbdd5fc8862a5!8m2!" },
    @{ Row = 8;  Name = "Springy Point Boat Launch";   Url = "https://www.google.com/maps/place/Springy+Point+Campground/@48.2362164,-116.5896929,797m/data=!3m1!1e3!4m9!1m2!2m1!1sspringy+point+boat+ramp!3m5!1s0x5363d435770b5855:0xb7b27f55e5104eca!8m2!3d48.2365207!4d-116.5861048!15sChdzcHJpbmd5IHBvaW50IGJvYXQgcmFtcFo" },
    @{ Row = 9;  Name = "Morton Slough Boat Launch";   Url = "https://www.google.com/maps/place/Morton+Slough+Boat+Launch/@48.1804787,-116.7319031,4321m/data=!3m1!1e3!4m9!1m2!2m1!1smorton+slough+boat+launch!3m5!1s0x53622a675485d28d:0x9ec2cc1934bd00a8!8m2!3d48.1804787!4d-116.7143936!15sChltb3J0b24gc2xvdWdoIGJvYXQgbGF" },
    @{ Row = 10; Name = "Riley Creek Boat Launch";     Url = "https://www.google.com/maps/place/Riley+Creek+Recreation+Area/@48.1598695,-116.7755195,510m/data=!3m1!1e3!4m9!1m2!2m1!1sriley+creek+boat+launch!3m5!1s0x536229ee6ad2d3c7:0x146b7d9edc002dab!8m2!3d48.1597415!4d-116.7744035!15sChdyaWxleSBjcmVlayBib2F0IGxhdW5j" },
    @{ Row = 11; Name = "Priest River Boat Launch";    Url = "https://www.google.com/maps/place/Bonner+Park+West/@48.1778442,-116.9080716,922m/data=!3m1!1e3!4m9!1m2!2m1!1spriest+river+boat+launch!3m5!1s0x0:0x587a07d691315514!8m2!3d48.1772667!4d-116.904774!15sChhwcmllc3Qgcml2ZXIgYm9hdCBsYXVuY2haJwoLYm9hdCBsYXVuY2giGH" },
    @{ Row = 12; Name = "Bayview Boat Launch";         Url = "https://www.google.com/maps/place/47%C2%B058'50.6%22N+116%C2%B033'30.4%22W/@47.9807186,-116.5606153,542m/data=!3m2!1e3!4b1!4m5!3m4!1s0x0:0x0!8m2!3d47.980715!4d-116.558432" },
    @{ Row = 13; Name = "Garfield Bay Boat Launch";    Url = "https://www.google.com/maps/place/Public+Boat+Launch/@48.1866988,-116.4432994,1361m/data=!3m1!1e3!4m9!1m2!2m1!1sgarfield+bay+boat+launch!3m5!1s0x0:0x73f46aca92e216f5!8m2!3d48.1878011!4d-116.4374398!15sChhnYXJmaWVsZCBiYXkgYm9hdCBsYXVuY2haJwoLYm9hdCBsYXVuY2" },
    @{ Row = 14; Name = "Johnson Creek Boat Launch";   Url = "https://www.google.com/maps/place/48%C2%B008'20.5%22N+116%C2%B013'43.4%22W/@48.1390409,-116.2299025,541m/data=!3m2!1e3!4b1!4m5!3m4!1s0x0:0x0!8m2!3d48.139039!4d-116.228712" },
    @{ Row = 15; Name = "Drift Yard Boat Launch";      Url = "https://www.google.com/maps/place/48%C2%B010'24.1%22N+116%C2%B013'55.5%22W/@48.1733514,-116.2329207,270m/data=!3m2!1e3!4b1!4m5!3m4!1s0x0:0x0!8m2!3d48.17335!4d-116.232091" },
    @{ Row = 16; Name = "Hope Boat Basin Boat Launch"; Url = "https://www.google.com/maps/place/Hope+Boat+Ramp/@48.2503613,-116.3161855,175m/data=!3m1!1e3!4m13!1m7!3m6!1s0x0:0x0!2zNDjCsDE1JzAxLjQiTiAxMTbCsDE4JzU2LjAiVw!3b1!8m2!3d48.250381!4d-116.315547!3m4!1s0x5361620ef5c0b96b:0xb95ad8326c224d3b!8m2!3d48.2504676!4d-" },
)

# 1. Move the "Info" hyperlink cell (and its style/format) from column E to
#    column F, preserving its text/format via copy+paste special so the
#    existing "Hyperlink" cell style is reused instead of a new one created.
for ($r = 2; $r -le 16; $r++) {
    $src = $ws.Range("E" + $r)
    $dst = $ws.Range("F" + $r)
    $src.Copy() | Out-Null
    $dst.PasteSpecial(-4163) | Out-Null   # xlPasteValues
    $src.Copy() | Out-Null
    $dst.PasteSpecial(-4122) | Out-Null   # xlPasteFormats
}
$excel.CutCopyMode = 0

# 2. Remove the old hyperlinks (column E) and add the new ones pointing at
#    column F, reusing the very same target URLs.
$ws.Range("E2:E16").Hyperlinks.Delete()
foreach ($s in $sites) {
    $ws.Hyperlinks.Add($ws.Range("F" + $s.Row), $s.Url, "", "", "click here") | Out-Null
}

# 3. Clear out the now-empty column E entirely (header + data) so the sheet
#    no longer has an "E" column at all.
$ws.Range("E1:E16").Clear() | Out-Null

# 4. Update the site name (column A) and boat-launch label (column D) on
#    every data row, and fill in the new "Size" column (G) with 50.
foreach ($s in $sites) {
    $r = $s.Row
    $ws.Range("A" + $r).Value = $s.Name
    $ws.Range("D" + $r).Value = "Boat Launch"
    $ws.Range("G" + $r).Value = 50
}

# 5. Header row: F1 becomes "Info" (moved from old E1) and G1 becomes "Size".
$ws.Range("F1").Value = "Info"
$ws.Range("G1").Value = "Size"

# 6. Misc cosmetic bits that came along with the re-save.
$ws.PageSetup.Orientation = 1
$ws.Range("I10").Select() | Out-Null

Write-Host "done"
